$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3
$ws.Range("O2").Value = 1
$ws.Range("X2").Value = "Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y2").Value = "BankStability_score,ChannelStability_score,Stability_Mean,CoarseSubstrate_score"
$ws.Range("Z2").Value = "BankStability_score,ChannelStability_score,Stability_Mean,CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J3").Value = 3
$ws.Range("O3").Value = 5
$ws.Range("Y3").Value = "CoarseSubstrate_score,Cover-Wood_score"
$ws.Range("Z3").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O4").Value = 1
$ws.Range("X4").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z4").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J5").Value = 3
$ws.Range("O5").Value = 5
$ws.Range("Y5").Value = "CoarseSubstrate_score"
$ws.Range("Z5").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J6").Value = 3
$ws.Range("O6").Value = 5
$ws.Range("Y6").Value = "CoarseSubstrate_score,Cover-Wood_score"
$ws.Range("Z6").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J7").Value = 3
$ws.Range("O7").Value = 3
$ws.Range("Y7").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z7").Value = "BankStability_score,ChannelStability_score,Stability_Mean,CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,Off-Channel-Floodplain_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J8").Value = 3
$ws.Range("O8").Value = 3
$ws.Range("Y8").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z8").Value = "BankStability_score,ChannelStability_score,Stability_Mean,CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J9").Value = 3
$ws.Range("O9").Value = 1
$ws.Range("X9").Value = "Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y9").Value = "CoarseSubstrate_score"
$ws.Range("Z9").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J10").Value = 3
$ws.Range("O10").Value = 5
$ws.Range("Y10").Value = "CoarseSubstrate_score"
$ws.Range("Z10").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O11").Value = 1
$ws.Range("X11").Value = "Flow-SummerBaseFlow_score,Off-Channel-Floodplain_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z11").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,Off-Channel-Floodplain_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J12").Value = 3
$ws.Range("O12").Value = 1
$ws.Range("X12").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y12").Value = "CoarseSubstrate_score"
$ws.Range("Z12").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O13").Value = 5
$ws.Range("O14").Value = 1
$ws.Range("X14").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z14").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J15").Value = 3
$ws.Range("O15").Value = 3
$ws.Range("Y15").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z15").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J16").Value = 3
$ws.Range("O16").Value = 1
$ws.Range("X16").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y16").Value = "CoarseSubstrate_score"
$ws.Range("Z16").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J17").Value = 3
$ws.Range("O17").Value = 1
$ws.Range("X17").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y17").Value = "CoarseSubstrate_score"
$ws.Range("Z17").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("J18").Value = 3
$ws.Range("O18").Value = 3
$ws.Range("Y18").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z18").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O19").Value = 5
$ws.Range("J20").Value = 3
$ws.Range("O20").Value = 1
$ws.Range("X20").Value = "Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y20").Value = "CoarseSubstrate_score"
$ws.Range("Z20").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O21").Value = 5
$ws.Range("J22").Value = 1
$ws.Range("O22").Value = 1
$ws.Range("X22").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Y22").ClearContents()
$ws.Range("Z22").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O23").Value = 3
$ws.Range("Y23").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z23").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O24").Value = 5
$ws.Range("O25").Value = 3
$ws.Range("Y25").Value = "CoarseSubstrate_score,PoolQuantity&Quality_score"
$ws.Range("Z25").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O26").Value = 5
$ws.Range("O27").Value = 5
$ws.Range("O28").Value = 5
$ws.Range("O29").Value = 5
$ws.Range("O30").Value = 1
$ws.Range("X30").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z30").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O31").Value = 1
$ws.Range("X31").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z31").Value = "CoarseSubstrate_score,Cover-Wood_score,Flow-SummerBaseFlow_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O32").Value = 5
$ws.Range("O33").Value = 1
$ws.Range("X33").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,Off-Channel-Floodplain_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("Z33").Value = "CoarseSubstrate_score,Flow-SummerBaseFlow_score,Off-Channel-Floodplain_score,PoolQuantity&Quality_score,Riparian-Disturbance_score,Riparian_Mean,Temperature-Rearing_score"
$ws.Range("O34").Value = 5
